$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param(
        [string]$addr,
        [string]$val
    )
    $r = $ws.Range($addr)
    $r.Value = "'" + $val
    $r.ClearFormats()
}

# Row 124
Set-TextCell "A124" "Murodullayev Asliddinjon Xurshidovich"
Set-TextCell "B124" "Yurisprudensiya"
Set-TextCell "C124" "O'zbek tili"
Set-TextCell "D124" "Kunduzgi"
Set-TextCell "E124" "AD4081875"
Set-TextCell "F124" "51705075820037"
Set-TextCell "G124" "Navoiy viloyati"
Set-TextCell "H124" "Karmana tumani"
Set-TextCell "I124" "998958222222"
Set-TextCell "J124" "+998958222222"
Set-TextCell "K124" "2025-06-27"

# Row 125
Set-TextCell "A125" "Azizbek Isroilov Xusniddinovich"
Set-TextCell "B125" "Yurisprudensiya"
Set-TextCell "C125" "O'zbek tili"
Set-TextCell "D125" "Kunduzgi"
Set-TextCell "E125" "AD6928826"
Set-TextCell "F125" "50201086580012"
Set-TextCell "G125" "Toshkent shahri"
Set-TextCell "H125" "Uchtepa tumani"
Set-TextCell "I125" "998909658251"
Set-TextCell "J125" "+998909658251"
Set-TextCell "K125" "2025-06-28"

# Row 126
Set-TextCell "A126" "orifjonov dilshod"
Set-TextCell "B126" "Yurisprudensiya"
Set-TextCell "C126" "O'zbek tili"
Set-TextCell "D126" "Kunduzgi"
Set-TextCell "E126" "AD1569162"
Set-TextCell "F126" "50509076920045"
Set-TextCell "G126" "Fargona viloyati"
Set-TextCell "H126" "Qo‘qon shahri"
Set-TextCell "I126" "998911398039"
Set-TextCell "J126" "+998911398039"
Set-TextCell "K126" "2025-06-28"

# Row 127
Set-TextCell "A127" "Xomidjonova Tamannoxon Axtamjon qizi"
Set-TextCell "B127" "Yurisprudensiya"
Set-TextCell "C127" "O'zbek tili"
Set-TextCell "D127" "Kunduzgi"
Set-TextCell "E127" "AD5616436"
Set-TextCell "F127" "62305077010104"
Set-TextCell "G127" "Fargona viloyati"
Set-TextCell "H127" "Fargʻona tumani"
Set-TextCell "I127" "998887602305"
Set-TextCell "J127" "+998507142104"
Set-TextCell "K127" "2025-06-28"
